$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column BO (LogReg output sigmut+others) values for rows 2-76
$boValues = @{
    2 = 0.707798182964325
    3 = 0.7345734238624573
    4 = 0.8727198243141174
    5 = 0.8167624473571777
    6 = 0.7034626007080078
    7 = 0.6226551532745361
    8 = 0.7518402934074402
    9 = 0.7594751715660095
    10 = 0.7373155355453491
    11 = 0.7546300292015076
    12 = 0.6847510933876038
    13 = 0.7531487345695496
    14 = 0.650474488735199
    15 = 0.6719115972518921
    16 = 0.7176676392555237
    17 = 0.7194453477859497
    18 = 0.7973012924194336
    19 = 0.4688084423542023
    20 = 0.7541728615760803
    21 = 0.4263558387756348
    22 = 0.6458883285522461
    23 = 0.8448514938354492
    24 = 0.8312001824378967
    25 = 0.7239437103271484
    26 = 0.7627598643302917
    27 = 0.834860622882843
    28 = 0.7576714158058167
    29 = 0.7069045901298523
    30 = 0.7657533884048462
    31 = 0.7477073073387146
    32 = 0.5232943892478943
    33 = 0.6623345613479614
    34 = 0.766302227973938
    35 = 0.7689911723136902
    36 = 0.5299518704414368
    37 = 0.7659668326377869
    38 = 0.7795442342758179
    39 = 0.466173529624939
    40 = 0.7568145394325256
    41 = 0.5414235591888428
    42 = 0.6749743819236755
    43 = 0.6900009512901306
    44 = 0.7264436483383179
    45 = 0.4763390123844147
    46 = 0.4493989646434784
    47 = 0.8399074673652649
    48 = 0.8476596474647522
    49 = 0.801840603351593
    50 = 0.8449442386627197
    51 = 0.7922573685646057
    52 = 0.5533796548843384
    53 = 0.1368321627378464
    54 = 0.07053074240684509
    55 = 0.7324074506759644
    56 = 0.7779104113578796
    57 = 0.8202462792396545
    58 = 0.9326225519180298
    59 = 0.7098234295845032
    60 = 0.8557806015014648
    61 = 0.8278502821922302
    62 = 0.8157700896263123
    63 = 0.6926621794700623
    64 = 0.7022183537483215
    65 = 0.7145379185676575
    66 = 0.1714053452014923
    67 = 0.6154088377952576
    68 = 0.3205895125865936
    69 = 0.8269971013069153
    70 = 0.8011592626571655
    71 = 0.6300471425056458
    72 = 0.6852648854255676
    73 = 0.8271722197532654
    74 = 0.7276062965393066
    75 = 0.8522377014160156
    76 = 0.8691616058349609
}

foreach ($row in $boValues.Keys) {
    $ws.Cells.Item($row, 67).Value = $boValues[$row]
}
